# edit.ps1 - applies the two changes described by the diff:
#   1. The table on the "B1- TYPES OF FINANCIAL DOCUMENTS" slide gets a new
#      table style ({C7759F40-...} -> {9F0E2676-...}).
#   2. The presentation's active theme colour scheme changes from the
#      "Red Violet" / Integral palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Update the table's style -------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{9F0E2676-700B-4BCA-9504-99ED4AE5F629}")
        }
    }
}

# --- 2. Swap the theme colour scheme back to the standard Office colours --------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # Dark 1     -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # Light 1    -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # Dark 2     -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # Light 2    -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # Accent 1   -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # Accent 2   -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # Accent 3   -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # Accent 4   -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # Accent 5   -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # Accent 6   -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # Hyperlink  -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # Followed Hyperlink -> 954F72
